# ElectroCalcs update — "Spinney and Matt updates and LHS implemented."
#
# motorCalcs: update input torque value.
# batteryCalcs: update input capacity value, add the LHS (least-squares /
#   regression) coefficient tables in columns J:L with labels in column K,
#   and re-point the weight/cost output formulas at the new quadratic /
#   power-law fits instead of the old two-point averages.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# motorCalcs
# ---------------------------------------------------------------------
$motor = $wb.Worksheets.Item("motorCalcs")
$motor.Range("B2").Value = 1.002

# ---------------------------------------------------------------------
# batteryCalcs
# ---------------------------------------------------------------------
$batt = $wb.Worksheets.Item("batteryCalcs")
$batt.Range("B2").Value = 11.989

# --- BatteryCost quadratic fit block (rows 1, 3:5, cols J:L) ---------
$batt.Range("K1").Value = "BatteryCost"
$batt.Range("L1").Formula = "=(J3*B2^2)+(J4*B2)+J5"

$batt.Range("J3").Value = 6.6303
$batt.Range("K3").Value = "6.6303x2 - 31.91x + 170.05"
$chars = $batt.Range("K3").Characters(8, 1)
$chars.Font.Superscript = $true

$batt.Range("J4").Value = -31.91
$batt.Range("J5").Value = 170.05

# --- BatteryWeight power fit block (row 8, J9:K9, J10:J11) -----------
$batt.Range("K8").Value = "BatteryWeight"
$batt.Range("L8").Formula = "=J14*B2^J15"

$batt.Range("J9").Value = -0.0355
$batt.Range("K9").Value = "y = -0.0355x2 + 0.7489x + 0.3598"
$chars2 = $batt.Range("K9").Characters(13, 1)
$chars2.Font.Superscript = $true

$batt.Range("J10").Value = 0.7489
$batt.Range("J11").Value = 0.3598

# --- second BatteryWeight power fit block (rows 13:19, J:K) ----------
# leading "=" would otherwise be auto-parsed as a formula; the classic
# text-qualifying apostrophe via .Formula forces literal text instead.
$batt.Range("K13").Formula = "'= 1.2097x0.5862"
$chars3 = $batt.Range("K13").Characters(10, 6)
$chars3.Font.Superscript = $true

$batt.Range("J14").Value = 1.2097
$batt.Range("J15").Value = 0.5862

$batt.Range("K16").Value = "BatteryWeight"
$batt.Range("L16").Formula = "=J18*B2^J19"

$batt.Range("K17").Value = "0.3424x + 1.3273"
$batt.Range("J18").Value = 0.3424
$batt.Range("J19").Value = 1.3273

# --- re-point the output formulas at the new fitted curves -----------
$batt.Range("E9").Formula = "=L16/2.2"
$batt.Range("E10").Formula = "=L1"

# ---------------------------------------------------------------------
# Sheet views — batteryCalcs becomes the active/selected sheet, with
# motorCalcs' selection left on F9 and batteryCalcs' on E9.
# ---------------------------------------------------------------------
$motor.Range("F9").Select()
$batt.Activate()
$batt.Range("E9").Select()
